$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 55
$ws.Range("I2").Value = 157
$ws.Range("J2").Value = 546
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 154
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 114
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 71
$ws.Range("T2").Value = 93
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 891
$ws.Range("X2").Value = 874
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 24
$ws.Range("AA2").Value = 8
